# Apply updated "dSF" (column F) values for several rows on Sheet1.
# These values were repulled/recalculated (per commit message:
# "repull data, push all data, mean calculation").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    9  = -6
    10 = -1
    11 = 0
    13 = 5
    15 = -3
    16 = 3
    17 = -2
    24 = 6
    26 = 10
    31 = 1
    32 = 2
    33 = -2
    34 = -1
    38 = 3
    44 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
